# daily auto push: 2026-02-04 22:45 UTC
# Insert a new data row just above the "2026/12/29" block (current row 773)
# for the 2026/02/05 (Thursday) 05:00 entry, shifting every following row
# down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push row 773 (and everything below it) down by one row.
$ws.Rows("773").Insert()

# Format the date cell as text first so the "YYYY/MM/DD" string is kept
# literally instead of being auto-converted into a date serial number
# (matching how every other date cell in column A is stored).
$ws.Range("A773").NumberFormat = "@"
$ws.Range("A773").Value = "2026/02/05"
$ws.Range("B773").Value = "木"
$ws.Range("C773").Value = 5
$ws.Range("D773").Value = 201
